$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    "C2" = 43;
    "C3" = 71;
    "C4" = 55;
    "C5" = 39;
    "C6" = 57;
    "C7" = 58;
    "C8" = 85;
    "C9" = 31;
    "C10" = 41;
    "C11" = 59;
    "C12" = 98;
    "C13" = 111;
    "C14" = 70;
    "C15" = 53;
    "C16" = 99;
    "C17" = 102;
    "C18" = 72;
    "C19" = 113;
    "C20" = 60;
    "C21" = 19;
    "C23" = 125;
    "C24" = 46;
    "C25" = 83;
    "C26" = 27;
    "C27" = 28;
    "C28" = 32;
    "C29" = 69;
    "C30" = 24;
    "C31" = 68;
    "C32" = 56;
    "C33" = 40;
    "C34" = 88;
    "C35" = 1;
    "C36" = 91;
    "C37" = 45;
    "C38" = 126;
    "C39" = 82;
    "C40" = 44;
    "C41" = 30;
    "C42" = 7
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
